$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Clear the numbering values (2..50) previously stored in B6:B54,
# leaving the cells blank (style is preserved automatically).
$ws.Range("B6:B54").ClearContents()

# Update the active selection to reflect the cleared range.
$ws.Range("B6:B54").Select()
